$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 152.22223
$ws.Range("I9").Value = 158.75
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 158.75
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = 10.25
$ws.Range("N9").Value = -438

# Row 74
$ws.Range("H74").Value = 6187.8125
$ws.Range("I74").Value = 6926.3335
$ws.Range("J74").Value = 3972.25
$ws.Range("K74").Value = 6926.3335
$ws.Range("L74").Value = 3972.25
$ws.Range("M74").Value = -5990.3335
$ws.Range("N74").Value = -5844.25

# Row 76
$ws.Range("H76").Value = 3592.5386
$ws.Range("I76").Value = 3643.2856
$ws.Range("J76").Value = 3533.3333
$ws.Range("K76").Value = 3643.2856
$ws.Range("L76").Value = 3533.3333
$ws.Range("M76").Value = -3328.2856
$ws.Range("N76").Value = -4163.3333

# Row 77
$ws.Range("H77").Value = 6187.8125
$ws.Range("I77").Value = 6926.3335
$ws.Range("J77").Value = 3972.25
$ws.Range("K77").Value = 34631.6675
$ws.Range("L77").Value = 19861.25
$ws.Range("M77").Value = -29951.6675
$ws.Range("N77").Value = -29221.25

# Row 79
$ws.Range("H79").Value = 3592.5386
$ws.Range("I79").Value = 3643.2856
$ws.Range("J79").Value = 3533.3333
$ws.Range("K79").Value = 3643.2856
$ws.Range("L79").Value = 3533.3333
$ws.Range("M79").Value = -2551.2856
$ws.Range("N79").Value = -5717.3333

# Row 86
$ws.Range("H86").Value = 2435.24
$ws.Range("I86").Value = 2115.2
$ws.Range("J86").Value = 2915.3
$ws.Range("K86").Value = 2115.2
$ws.Range("L86").Value = 2915.3
$ws.Range("M86").Value = -992.1999999999998
$ws.Range("N86").Value = -5161.3

# Row 89
$ws.Range("H89").Value = 2435.24
$ws.Range("I89").Value = 2115.2
$ws.Range("J89").Value = 2915.3
$ws.Range("K89").Value = 10576
$ws.Range("L89").Value = 14576.5
$ws.Range("M89").Value = -4960
$ws.Range("N89").Value = -25808.5

# Row 132
$ws.Range("H132").Value = 6948633.5
$ws.Range("I132").Value = 9805968
$ws.Range("J132").Value = 9392.857
$ws.Range("K132").Value = 29417904
$ws.Range("L132").Value = 28178.571
$ws.Range("M132").Value = -29415374
$ws.Range("N132").Value = -33238.571

# Row 138
$ws.Range("H138").Value = 1531.65
$ws.Range("I138").Value = 704.9429
$ws.Range("J138").Value = 1976.8
$ws.Range("K138").Value = 2114.8287
$ws.Range("L138").Value = 5930.4
$ws.Range("M138").Value = 3025.1713
$ws.Range("N138").Value = -16210.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3901.05
$ws.Range("I32").Value = 4336
$ws.Range("K32").Value = 4336
$ws.Range("M32").Value = -4049

# Row 132
$ws.Range("H132").Value = 3514.5
$ws.Range("I132").Value = 3221.3635
$ws.Range("K132").Value = 9664.0905
$ws.Range("M132").Value = -7134.0905

# Row 133
$ws.Range("H133").Value = 30915
$ws.Range("J133").Value = 30915
$ws.Range("L133").Value = 30915
$ws.Range("N133").Value = -35975

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1119.8462
$ws.Range("I31").Value = 1013.1667
$ws.Range("K31").Value = 1013.1667
$ws.Range("M31").Value = -718.1667

# Row 34
$ws.Range("H34").Value = 1119.8462
$ws.Range("I34").Value = 1013.1667
$ws.Range("K34").Value = 1013.1667
$ws.Range("M34").Value = -811.1667

# Row 109
$ws.Range("H109").Value = 15500.5
$ws.Range("J109").Value = 15500.5
$ws.Range("L109").Value = 15500.5
$ws.Range("N109").Value = -17580.5

# Row 132
$ws.Range("H132").Value = 9190.5
$ws.Range("I132").Value = 14975.125
$ws.Range("J132").Value = 3405.875
$ws.Range("K132").Value = 44925.375
$ws.Range("L132").Value = 10217.625
$ws.Range("M132").Value = -42395.375
$ws.Range("N132").Value = -15277.625

# Row 134
$ws.Range("H134").Value = 15153186
$ws.Range("I134").Value = 18520072
$ws.Range("J134").Value = 2200
$ws.Range("K134").Value = 55560216
$ws.Range("L134").Value = 6600
$ws.Range("M134").Value = -55557681
$ws.Range("N134").Value = -11670

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 146.45454
$ws.Range("J2").Value = 262.4
$ws.Range("L2").Value = 1574.4
$ws.Range("N2").Value = -1800.4

# Row 113
$ws.Range("H113").Value = 760.5
$ws.Range("I113").Value = 897.6667
$ws.Range("J113").Value = 738.8421
$ws.Range("K113").Value = 2693.0001
$ws.Range("L113").Value = 2216.5263
$ws.Range("M113").Value = -523.0001000000002
$ws.Range("N113").Value = -6556.5263

# Row 131
$ws.Range("H131").Value = 15385907
$ws.Range("I131").Value = 125000390
$ws.Range("J131").Value = 1418.0526
$ws.Range("K131").Value = 375001170
$ws.Range("L131").Value = 4254.1578
$ws.Range("M131").Value = -374996130
$ws.Range("N131").Value = -14334.1578

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 7067
$ws.Range("J80").Value = 7067
$ws.Range("L80").Value = 7067
$ws.Range("N80").Value = -9063

# Row 83
$ws.Range("H83").Value = 7067
$ws.Range("J83").Value = 7067
$ws.Range("L83").Value = 35335
$ws.Range("N83").Value = -45319

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 47830.816
$ws.Range("I132").Value = 2388.6365
$ws.Range("J132").Value = 93273
$ws.Range("K132").Value = 7165.9095
$ws.Range("L132").Value = 279819
$ws.Range("M132").Value = -4635.9095
$ws.Range("N132").Value = -284879

# Row 136
$ws.Range("H136").Value = 7211.722
$ws.Range("I136").Value = 11330.8
$ws.Range("K136").Value = 33992.39999999999
$ws.Range("M136").Value = -31442.39999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 2400
$ws.Range("I15").Value = 1563.6364
$ws.Range("K15").Value = 1563.6364
$ws.Range("M15").Value = -1275.6364

# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Row 68
$ws.Range("H68").Value = 75000
$ws.Range("J68").Value = 75000
$ws.Range("L68").Value = 75000
$ws.Range("N68").Value = -76622

# Row 71
$ws.Range("H71").Value = 75000
$ws.Range("J71").Value = 75000
$ws.Range("L71").Value = 225000
$ws.Range("N71").Value = -233112

# Row 81
$ws.Range("H81").Value = 4620.3
$ws.Range("J81").Value = 5934.3335
$ws.Range("L81").Value = 11868.667
$ws.Range("N81").Value = -13990.667

# Row 84
$ws.Range("H84").Value = 4620.3
$ws.Range("J84").Value = 5934.3335
$ws.Range("L84").Value = 59343.335
$ws.Range("N84").Value = -69951.33499999999

# Row 132
$ws.Range("H132").Value = 3645.7
$ws.Range("I132").Value = 2994.4666
$ws.Range("K132").Value = 8983.399800000001
$ws.Range("M132").Value = -6453.399800000001

# Row 136
$ws.Range("H136").Value = 835
$ws.Range("I136").Value = 472.375
$ws.Range("K136").Value = 1417.125
$ws.Range("M136").Value = 1132.875
